$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.199161529541016
$ws.Range("B1").Value = 2.014509439468384
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.036058187484741
$ws.Range("E1").Value = 1.208304166793823
